$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of existing value in B2 ("supermercado" -> "Supermercado")
$ws.Range("B2").Value = "Supermercado"

# Add new row 3, copying the format of A2 into A3 (border/alignment style)
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Mayorista"
